$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the nine separate title runs ("Applications" " " "of" " " ...)
#    into a single run with the full title text, matching the target OOXML.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Applications of Gaussian Elimination Questions",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Applications of Gaussian Elimination Questions", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) For every equation (m:oMath) in the document, swap the order of the
#    <m:endChr/> and <m:sepChr/> children of <m:dPr/> so that <m:sepChr/>
#    comes immediately after <m:begChr/> (and before <m:endChr/>).
#    There is no dedicated OMathDelim.SepChar/EndChar property on the Word
#    object model, so we round-trip each equation's XML via
#    Range.WordOpenXML / Range.InsertXML, rewriting only the element order.
# ---------------------------------------------------------------------------
$count = $d.OMaths.Count

for ($i = 1; $i -le $count; $i++) {
    $om = $d.OMaths($i)
    $full = $om.Range.WordOpenXML

    # Pull out the <w:document>...</w:document> payload for this range.
    $partMarker = '<pkg:part pkg:name="/word/document.xml"'
    $partIdx = $full.IndexOf($partMarker)
    $xmlDataIdx = $full.IndexOf("<pkg:xmlData>", $partIdx)
    $docStart = $xmlDataIdx + "<pkg:xmlData>".Length
    $docEnd = $full.IndexOf("</pkg:xmlData></pkg:part>", $docStart)
    $docXml = $full.Substring($docStart, $docEnd - $docStart)

    # Locate the <m:oMath>...</m:oMath> fragment.
    $mathStartTag = "<m:oMath>"
    $mathEndTag = "</m:oMath>"
    $mathStart = $docXml.IndexOf($mathStartTag)
    $mathEnd = $docXml.IndexOf($mathEndTag, $mathStart) + $mathEndTag.Length

    # The equation is normally wrapped in <m:oMathPara>[<m:oMathParaPr/>]...
    # Include that wrapper in the replacement payload so InsertXML doesn't
    # drop it (InsertXML replaces exactly what the range used to contain).
    $paraTag = "<m:oMathPara>"
    $fragStart = $mathStart
    $rootIsPara = $false
    $beforeMath = $docXml.Substring(0, $mathStart)

    if ($beforeMath.EndsWith($paraTag)) {
        $rootIsPara = $true
        $fragStart = $mathStart - $paraTag.Length
    } else {
        $paraPrEndTag = "</m:oMathParaPr>"
        $paraPrEndIdx = $beforeMath.LastIndexOf($paraPrEndTag)
        if ($paraPrEndIdx -ge 0 -and ($paraPrEndIdx + $paraPrEndTag.Length) -eq $beforeMath.Length) {
            $paraStartIdx = $beforeMath.LastIndexOf($paraTag)
            if ($paraStartIdx -ge 0) {
                $rootIsPara = $true
                $fragStart = $paraStartIdx
            }
        }
    }

    $fragEnd = $mathEnd
    if ($rootIsPara) {
        $paraEndTag = "</m:oMathPara>"
        $fragEnd = $docXml.IndexOf($paraEndTag, $mathEnd) + $paraEndTag.Length
    }

    $fragXml = $docXml.Substring($fragStart, $fragEnd - $fragStart)

    # Swap <m:endChr .../><m:sepChr .../>  ->  <m:sepChr .../><m:endChr .../>
    $swapPattern = '(<m:endChr m:val="[^"]*"\s*/>)(<m:sepChr m:val="[^"]*"\s*/>)'
    $newFragXml = [System.Text.RegularExpressions.Regex]::Replace($fragXml, $swapPattern, '$2$1')

    if ($newFragXml -ne $fragXml) {
        $om.Range.InsertXML($newFragXml)
    }
}
